# Apply updated dSF ("F") column values for the specified rows.
# Mapping of worksheet row number -> new value for column F.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -8
    7  = -2
    10 = -1
    11 = -3
    13 = -5
    16 = 0
    17 = -2
    22 = -6
    28 = 2
    29 = 2
    32 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
